# Update the "Skills" list: row 13 (B13) now reads "Trade floor support role5."
# instead of "Trade floor support role". All the other rows keep their text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B13").Value = "Trade floor support role5."

# Leave the selection on the edited cell, matching the saved view state.
$ws.Range("B13").Select()
